# Automatische test-sync: 2025-06-22 18:51:50
#
# Adds a new "Verzoek om factuur" row to the Logs sheet, introduces a new
# "Factuur / Administratie" category (reshuffling the Dashboard summary
# rows to match), and extends the conditional formatting / chart ranges
# that depended on the old row counts.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new incoming mail as row 21
# ---------------------------------------------------------------------
$logs.Range("A21").Value = "Verzoek om factuur"
$logs.Range("B21").Value = "mailmind.test@zohomail.eu"
$logs.Range("C21").Value = "Kunt u mij een factuur sturen voor mijn laatste bestelling?"
$logs.Range("D21").Value = "Factuur / Administratie"
$logs.Range("E21").Value = "Beste klant,`nBedankt voor uw e-mail. Om u te kunnen helpen met het versturen van een factuur voor uw laatste bestelling, hebben wij wat extra informatie nodig. Kunt u ons alstublieft de volgende gegevens verstrekken:`n1. Uw bestelnummer`n2. De datum van uw bestelling`n3. Het e-mailadres waarnaar wij de factuur kunnen sturen`nZodra wij deze gegevens van u hebben ontvangen, zullen wij zo spoedig mogelijk de factuur voor u opstellen en toesturen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F21").Value = "2025-06-22 18:51:15"
$logs.Range("G21").Value = "Ja"

# Extend the two conditional-formatting blocks (category / answered)
# so they keep covering the data range now that row 21 exists.
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: re-sync the per-category counts.
#    The new "Factuur / Administratie" category (count 1) slots in among
#    the other count==1 categories, nudging the remaining ties around;
#    rewrite the affected rows (4-12) in their new order and append the
#    new row 13.
# ---------------------------------------------------------------------
$dash.Range("A4").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B4").Value = 2

$dash.Range("A5").Value = "Sollicitatie / Vacature"
$dash.Range("B5").Value = 2

$dash.Range("A6").Value = "Overig"
$dash.Range("B6").Value = 2

$dash.Range("A9").Value = "Openingstijden / Locatie"
$dash.Range("B9").Value = 1

$dash.Range("A10").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B10").Value = 1

$dash.Range("A11").Value = "Offerte / Prijsaanvraag"
$dash.Range("B11").Value = 1

$dash.Range("A12").Value = "Retour / Terugbetaling"
$dash.Range("B12").Value = 1

$dash.Range("A13").Value = "Factuur / Administratie"
$dash.Range("B13").Value = 1

# ---------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the category/value series
#    references to cover the new row 13.
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$13,Dashboard!`$B`$2:`$B`$13,1)"
